$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), matching the style of the other header cells (copy G1's
# formatting into H1 first so the new header cell reuses the existing bold/bordered header style)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the new Save column values for the two data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
